# Update import/export in hr and setting:
# - Remove the "id" column (column A), shifting the remaining
#   name/email/phone/address/status/group_id columns one to the left.
# - Remove the leftover duplicated "Jan 19, 2025" columns (old H and I,
#   which become G and H after the id column is removed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "id" column entirely.
$ws.Columns("A").Delete()

# Select and delete the two stray date columns that trailed the data
# (now at G and H after the shift above).
$ws.Columns("G:G").Select()
$ws.Columns("G:H").Delete()
